$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 410 - this shifts the existing row 410
# (and everything below it) down to row 411, matching the rest of the sheet.
$ws.Rows("410:410").Insert()

# Populate the newly inserted row 410 with the new data point.
$ws.Range("A410").Value = 11
$ws.Range("B410").Value = "Vega Monumental Concepción"
$ws.Range("C410").Value = "Bíobío"
$ws.Range("D410").Value = 45209
$ws.Range("E410").Value = 8
$ws.Range("F410").Value = 100112008
$ws.Range("G410").Value = "Coliflor"
$ws.Range("H410").Value = "Sin especificar"
$ws.Range("I410").Value = "Primera"
$ws.Range("J410").Value = 800
$ws.Range("K410").Value = 1000
$ws.Range("L410").Value = 1000
$ws.Range("M410").Value = 1000
$ws.Range("N410").Value = "$/unidad"
$ws.Range("O410").Value = "Región Metropolitana"
$ws.Range("P410").Value = 1000
$ws.Range("Q410").Value = 1
$ws.Range("R410").Value = "Hortaliza"
